# This workbook holds a weekly price log for "Bruselas (repollito)" at
# Vega Central Mapocho de Santiago. A new weekly record needs to be
# inserted at row 13 (pushing the existing row 13 and all rows below it
# down by one row), matching the structure of every other data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting rows 13:82 down to 14:83.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly entry.
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 45071
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100112035
$ws.Range("G13").Value = "Bruselas (repollito)"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 52
$ws.Range("K13").Value = 22000
$ws.Range("L13").Value = 24000
$ws.Range("M13").Value = 23000
$ws.Range("N13").Value = "$/malla 15 kilos"
$ws.Range("O13").Value = "Provincia de Quillota"
$ws.Range("P13").Value = 1533
$ws.Range("Q13").Value = 15
$ws.Range("R13").Value = "Hortaliza"
